# Regenerate save_data to use K (strikeouts) instead of Strike# for the
# pitching log on Sheet1. Column G ("K") values are recalculated for each
# outing row and rewritten in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row number -> new K value (column G)
$newK = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 1
    6  = 2
    7  = 2
    8  = 1
    9  = 1
    10 = 0
    11 = 1
    14 = 2
    15 = 2
    16 = 1
    17 = 2
    18 = 0
    19 = 1
    20 = 0
    21 = 1
    22 = 2
    23 = 2
    24 = 0
    25 = 0
    26 = 0
    27 = 2
    28 = 0
    29 = 1
    30 = 2
    31 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
